$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi output updated following Dr Hou advice:
# the "sCs" sending-cluster category is now included alongside the
# existing "ECs" and "FAPs" sending clusters for the Sfrp2-Fzd5 pair,
# and per-cluster ligand/receptor/edge statistics are recomputed.

$rows = @(
    @{ Row = 2; A = "ECs"; B = "Sfrp2"; C = "Fzd5"; D = "ECs"; Nums = @(3, 1, 6.575715666666667, 19.727147, 0.2077500326255795, 0.2077500326255795, 3, 1, 9.920975666666667, 29.762927, 0.5379109985456686, 0.5379109985456685, 65.23751511991878, 587.1376360792691, 0.1117510274975207, 0.1117510274975207) }
    @{ Row = 3; A = "ECs"; B = "Sfrp2"; C = "Fzd5"; D = "FAPs"; Nums = @(3, 1, 6.575715666666667, 19.727147, 0.2077500326255795, 0.2077500326255795, 3, 1, 6.153936333333334, 18.461809, 0.3336637594195427, 0.3336637594195427, 40.46653555876923, 364.1988200289231, 0.06931865690538348, 0.06931865690538348) }
    @{ Row = 4; A = "ECs"; B = "Sfrp2"; C = "Fzd5"; D = "sCs"; Nums = @(3, 1, 6.575715666666667, 19.727147, 0.2077500326255795, 0.2077500326255795, 3, 1, 2.368614333333333, 7.105843, 0.1284252420347888, 0.1284252420347888, 15.57533437999122, 140.178009419921, 0.0266803482226753, 0.0266803482226753) }
    @{ Row = 5; A = "FAPs"; B = "Sfrp2"; C = "Fzd5"; D = "ECs"; Nums = @(3, 1, 24.874192, 74.622576, 0.7858633890954824, 0.7858633890954825, 3, 1, 9.920975666666667, 29.762927, 0.5379109985456686, 0.5379109985456685, 246.7762535599946, 2220.986282039952, 0.4227245603488342, 0.4227245603488342) }
    @{ Row = 6; A = "FAPs"; B = "Sfrp2"; C = "Fzd5"; D = "FAPs"; Nums = @(3, 1, 24.874192, 74.622576, 0.7858633890954824, 0.7858633890954825, 3, 1, 6.153936333333334, 18.461809, 0.3336637594195427, 0.3336637594195427, 153.0741939111093, 1377.667745199984, 0.2622141327957815, 0.2622141327957815) }
    @{ Row = 7; A = "FAPs"; B = "Sfrp2"; C = "Fzd5"; D = "sCs"; Nums = @(3, 1, 24.874192, 74.622576, 0.7858633890954824, 0.7858633890954825, 3, 1, 2.368614333333333, 7.105843, 0.1284252420347888, 0.1284252420347888, 58.91736770128533, 530.256309311568, 0.1009246959508667, 0.1009246959508667) }
    @{ Row = 8; A = "sCs"; B = "Sfrp2"; C = "Fzd5"; D = "ECs"; Nums = @(1, 0.3333333333333333, 0.2021483333333333, 0.606445, 0.006386578278938132, 0.006386578278938132, 3, 1, 9.920975666666667, 29.762927, 0.5379109985456686, 0.5379109985456685, 2.005508696057222, 18.049578264515, 0.003435410699313688, 0.003435410699313687) }
    @{ Row = 9; A = "sCs"; B = "Sfrp2"; C = "Fzd5"; D = "FAPs"; Nums = @(1, 0.3333333333333333, 0.2021483333333333, 0.606445, 0.006386578278938132, 0.006386578278938132, 3, 1, 6.153936333333334, 18.461809, 0.3336637594195427, 0.3336637594195427, 1.244007973222778, 11.196071759005, 0.00213096971837769, 0.00213096971837769) }
    @{ Row = 10; A = "sCs"; B = "Sfrp2"; C = "Fzd5"; D = "sCs"; Nums = @(1, 0.3333333333333333, 0.2021483333333333, 0.606445, 0.006386578278938132, 0.006386578278938132, 3, 1, 2.368614333333333, 7.105843, 0.1284252420347888, 0.1284252420347888, 0.4788114397927778, 4.309302958135, 0.0008201978612467542, 0.0008201978612467542) }
)

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($r in $rows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    for ($i = 0; $i -lt $numCols.Length; $i++) {
        $ws.Range("$($numCols[$i])$($r.Row)").Value = $r.Nums[$i]
    }
}
